$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Highlight the q=0.25 "p" totals column (G9:G11) with a new green fill,
# matching the ad-hoc highlighting applied to q=0.75 (Q8:Q10) in the
# existing "totals" orange used elsewhere on the sheet (B10:B12, L10/L12/L13).
$ws.Range("G9:G11").Interior.Color = 5296274   # RGB(146,208,80) = 92D050 (green)
$ws.Range("Q8:Q10").Interior.Color = 49407     # RGB(255,192,0)  = FFC000 (orange)

# Leave the selection where the author last clicked.
$ws.Range("L19").Select()
